# Update faturamento_diario_lojas.xlsx: fill in column AC (value for a
# previously-zero day) for each store row, and recompute the row total
# in column AG (sum of columns B:AF) to reflect the new AC value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bibi Cell Mundi
$ws.Range("AC2").Value = 11700.63
$ws.Range("AG2").Value = 307481.48

# Row 3: Bibi Cell Vieiralves
$ws.Range("AC3").Value = 6882
$ws.Range("AG3").Value = 193807.3

# Row 4: Bibi Cell Ponta Negra
$ws.Range("AC4").Value = 1639
$ws.Range("AG4").Value = 93357.33

# Row 5: Bibi Cell Manauara
$ws.Range("AC5").Value = 1398
$ws.Range("AG5").Value = 82795.48

# Row 6: total (sum of rows 2-5)
$ws.Range("AC6").Value = 21619.63
$ws.Range("AG6").Value = 677441.59
